# Insert two new data rows at row 353 (pushing the existing rows 353:451
# down to 355:453) and populate the two new rows with their values.
# This mirrors the target diff: the dimension grows from A1:R451 to
# A1:R453, all former rows 353-451 become rows 355-453 unchanged, and the
# two newly inserted rows 353-354 carry new "Poroto granado" price entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above current row 353 (one at a time), which
# shifts the old rows 353-451 down to become rows 355-453.
$ws.Rows.Item(353).Insert()
$ws.Rows.Item(353).Insert()

# ---- New row 353 ----
$ws.Range("A353").Value = 6
$ws.Range("B353").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C353").Value = "Metropolitana"
$ws.Range("D353").Value = 44627
$ws.Range("E353").Value = 13
$ws.Range("F353").Value = 100112030
$ws.Range("G353").Value = "Poroto granado"
$ws.Range("H353").Value = "Sin especificar"
$ws.Range("I353").Value = "Primera"
$ws.Range("J353").Value = 300
$ws.Range("K353").Value = 18000
$ws.Range("L353").Value = 20000
$ws.Range("M353").Value = 18867
$ws.Range("N353").Value = "`$/saco 25 kilos"
$ws.Range("O353").Value = "Región Metropolitana"
$ws.Range("P353").Value = 755
$ws.Range("Q353").Value = 25
$ws.Range("R353").Value = "Hortaliza"

# ---- New row 354 ----
$ws.Range("A354").Value = 6
$ws.Range("B354").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C354").Value = "Metropolitana"
$ws.Range("D354").Value = 44627
$ws.Range("E354").Value = 13
$ws.Range("F354").Value = 100112030
$ws.Range("G354").Value = "Poroto granado"
$ws.Range("H354").Value = "Sin especificar"
$ws.Range("I354").Value = "Primera"
$ws.Range("J354").Value = 630
$ws.Range("K354").Value = 17000
$ws.Range("L354").Value = 20000
$ws.Range("M354").Value = 18190
$ws.Range("N354").Value = "`$/saco 25 kilos"
$ws.Range("O354").Value = "Región de O'Higgins"
$ws.Range("P354").Value = 728
$ws.Range("Q354").Value = 25
$ws.Range("R354").Value = "Hortaliza"
